$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "credit_supervised"
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = "python"
$ws.Range("E8").Value = "Credit Risk"
$ws.Range("F8").Value = "https://github.com/cdpeters/credit-risk-supervised-ML-sklearn"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://github.com/cdpeters/credit-risk-supervised-ML-sklearn")
$ws.Range("G8").Value = "python"
$ws.Range("H8").Value = "pandas, sklearn, imblearn"
$ws.Range("I8").Value = "resampling: oversampling, undersampling, combination sampling ensemble methods: random forest, AdaBoost"

# Copy formatting from row 6 into row 8 so styles match exactly (reused cell
# styles rather than new ones synthesized by Hyperlinks.Add)
$ws.Range("A6:I6").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)  # xlPasteFormats
